# Updates the LR-pairs table (Adm-Ramp2) with refreshed TPM-derived values.
# The underlying computation now also covers the Resolving-Mac sending
# cluster (previously missing), so the table grows from 12 to 16 data rows
# (rows 2-17) while the header row (row 1) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @('ECs','Adm','Ramp2','ECs',3,1,13.26223533333333,39.786706,0.391007655706778,0.3910076557067781,3,1,82.374849,247.124547,0.6248537741508955,0.6248537741508954,1092.474632985798,9832.271696872182,0.2443226093902741,0.2443226093902742),
    @('ECs','Adm','Ramp2','FAPs',3,1,13.26223533333333,39.786706,0.391007655706778,0.3910076557067781,3,1,32.15646633333333,96.46939900000001,0.2439226243891451,0.2439226243891451,426.466624001077,3838.199616009694,0.09537561353624456,0.09537561353624459),
    @('ECs','Adm','Ramp2','MuSCs',3,1,13.26223533333333,39.786706,0.391007655706778,0.3910076557067781,3,1,17.22395566666667,51.671867,0.130652181275918,0.1306521812759179,228.4281534222335,2055.853380800102,0.05108600311367367,0.05108600311367368),
    @('ECs','Adm','Ramp2','Resolving-Mac',3,1,13.26223533333333,39.786706,0.391007655706778,0.3910076557067781,2,0.6666666666666666,0.07533066666666667,0.225992,0.0005714201840414873,0.0005714201840414872,0.9990530291502221,8.991477262352,0.0002234296665855976,0.0002234296665855976),
    @('FAPs','Adm','Ramp2','ECs',3,1,19.55844,58.67532,0.5766373150128344,0.5766373150128344,3,1,82.374849,247.124547,0.6248537741508955,0.6248537741508954,1611.12354167556,14500.11187508004,0.3603140026020084,0.3603140026020084),
    @('FAPs','Adm','Ramp2','FAPs',3,1,19.55844,58.67532,0.5766373150128344,0.5766373150128344,3,1,32.15646633333333,96.46939900000001,0.2439226243891451,0.2439226243891451,628.93031739252,5660.37285653268,0.1406548871986407,0.1406548871986407),
    @('FAPs','Adm','Ramp2','MuSCs',3,1,19.55844,58.67532,0.5766373150128344,0.5766373150128344,3,1,17.22395566666667,51.671867,0.130652181275918,0.1306521812759179,336.87370346916,3031.86333122244,0.07533892301151546,0.07533892301151544),
    @('FAPs','Adm','Ramp2','Resolving-Mac',3,1,19.55844,58.67532,0.5766373150128344,0.5766373150128344,2,0.6666666666666666,0.07533066666666667,0.225992,0.0005714201840414873,0.0005714201840414872,1.47335032416,13.26015291744,0.000329502200669823,0.0003295022006698229),
    @('MuSCs','Adm','Ramp2','ECs',3,1,0.8919193333333334,2.675758,0.02629626747232247,0.02629626747232247,3,1,82.374849,247.124547,0.6248537741508955,0.6248537741508954,73.471720403514,661.245483631626,0.01643132197616212,0.01643132197616212),
    @('MuSCs','Adm','Ramp2','FAPs',3,1,0.8919193333333334,2.675758,0.02629626747232247,0.02629626747232247,3,1,32.15646633333333,96.46939900000001,0.2439226243891451,0.2439226243891451,28.68097401438245,258.1287661294421,0.006414254573487806,0.006414254573487807),
    @('MuSCs','Adm','Ramp2','MuSCs',3,1,0.8919193333333334,2.675758,0.02629626747232247,0.02629626747232247,3,1,17.22395566666667,51.671867,0.130652181275918,0.1306521812759179,15.36237905557622,138.261411500186,0.0034356647046739,0.0034356647046739),
    @('MuSCs','Adm','Ramp2','Resolving-Mac',3,1,0.8919193333333334,2.675758,0.02629626747232247,0.02629626747232247,2,0.6666666666666666,0.07533066666666667,0.225992,0.0005714201840414873,0.0005714201840414872,0.0671888779928889,0.604699901936,0.00001502621799863868,0.00001502621799863868),
    @('Resolving-Mac','Adm','Ramp2','ECs',2,0.6666666666666666,0.2055016666666667,0.6165050000000001,0.006058761808064916,0.006058761808064917,3,1,82.374849,247.124547,0.6248537741508955,0.6248537741508954,16.928168760915,152.353518848235,0.003785840182450667,0.003785840182450667),
    @('Resolving-Mac','Adm','Ramp2','FAPs',2,0.6666666666666666,0.2055016666666667,0.6165050000000001,0.006058761808064916,0.006058761808064917,3,1,32.15646633333333,96.46939900000001,0.2439226243891451,0.2439226243891451,6.608207425610557,59.47386683049501,0.001477869080771916,0.001477869080771916),
    @('Resolving-Mac','Adm','Ramp2','MuSCs',2,0.6666666666666666,0.2055016666666667,0.6165050000000001,0.006058761808064916,0.006058761808064917,3,1,17.22395566666667,51.671867,0.130652181275918,0.1306521812759179,3.539551596092778,31.855964364835,0.0007915904460549059,0.0007915904460549059),
    @('Resolving-Mac','Adm','Ramp2','Resolving-Mac',2,0.6666666666666666,0.2055016666666667,0.6165050000000001,0.006058761808064916,0.006058761808064917,2,0.6666666666666666,0.07533066666666667,0.225992,0.0005714201840414873,0.0005714201840414872,0.01548057755111111,0.13932519796,0.000003462098787427989,0.000003462098787427989)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $i + 2
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
